# Ebi3-Il27ra.xlsx: refresh LR-pair table with newly computed TPM values.
#
# The "ECs" sending-cluster rows (old rows 2-4) are dropped entirely, and the
# "Resolving-Mac" sending-cluster rows (old rows 5-7) move up to become the
# new rows 2-4, carrying freshly recomputed expression/specificity numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "ECs" sending-cluster rows; rows 5-7 ("Resolving-Mac") shift
# up to rows 2-4, keeping columns A:D (Sending/Ligand/Receptor/Target) intact.
$ws.Rows("2:4").Delete()

# New TPM-derived values for each (now-shifted) data row.
$rowData = @{
    2 = @{
        E = 3; F = 1; G = 9.105765666666665; H = 27.317297
        I = 1; J = 1
        K = 3; L = 1; M = 2.185458333333334; N = 6.556375
        O = 0.4976439168256567; P = 0.4976439168256567
        Q = 19.90027145759722; R = 179.102443118375
        S = 0.4976439168256567; T = 0.4976439168256567
    }
    3 = @{
        E = 3; F = 1; G = 9.105765666666665; H = 27.317297
        I = 1; J = 1
        K = 3; L = 1; M = 1.742815333333333; N = 5.228446
        O = 0.396851056620684; P = 0.396851056620684
        Q = 15.86966802560689; R = 142.827012230462
        S = 0.396851056620684; T = 0.396851056620684
    }
    4 = @{
        E = 3; F = 1; G = 9.105765666666665; H = 27.317297
        I = 1; J = 1
        K = 3; L = 1; M = 0.4633369999999999; N = 1.390011
        O = 0.1055050265536593; P = 0.1055050265536593
        Q = 4.219038146696332; R = 37.97134332026699
        S = 0.1055050265536593; T = 0.1055050265536593
    }
}

foreach ($r in $rowData.Keys) {
    foreach ($col in $rowData[$r].Keys) {
        $ws.Range("$col$r").Value2 = $rowData[$r][$col]
    }
}
